$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.17949914932251
$ws.Range("B1").Value = 2.415983200073242
$ws.Range("D1").Value = 2.333372354507446
$ws.Range("E1").Value = 1.194314479827881
